$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Cells whose new value is a bare percentage need NumberFormat forced to
# text first, otherwise Excel auto-converts "NN%" strings into a numeric
# percentage value instead of keeping literal text.
$percentCells = @("H4","H5","H6","H12","H14","H15","H21","H22","H30","H32","H33","H34")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("E2").Value = "2026-02-05 18:48:07"

# Row 3
$ws.Range("E3").Value = "2026-02-05 18:48:10"
$ws.Range("I3").Value = "9.6 mm"

# Row 4
$ws.Range("E4").Value = "2026-02-05 18:48:12"
$ws.Range("H4").Value = "68%"
$ws.Range("J4").Value = "990.0 hPa"
$ws.Range("O4").Value = "11.2 °C"

# Row 5
$ws.Range("E5").Value = "2026-02-05 18:48:15"
$ws.Range("H5").Value = "72%"
$ws.Range("J5").Value = "990.2 hPa"
$ws.Range("L5").Value = "45.7 km/h - 276º 18:22 TU"
$ws.Range("O5").Value = "9.6 °C"

# Row 6
$ws.Range("E6").Value = "2026-02-05 18:48:18"
$ws.Range("H6").Value = "73%"
$ws.Range("I6").Value = "1.0 mm"
$ws.Range("L6").Value = "36.4 km/h - 279º 18:18 TU"
$ws.Range("M6").Value = "16.5 °C 18:29 TU"
$ws.Range("O6").Value = "12.7 °C"

# Row 7
$ws.Range("E7").Value = "2026-02-05 18:48:20"
$ws.Range("I7").Value = "1.4 mm"
$ws.Range("J7").Value = "991.7 hPa"
$ws.Range("O7").Value = "10.3 °C"

# Row 8
$ws.Range("E8").Value = "2026-02-05 18:48:23"
$ws.Range("M8").Value = "14.6 °C 18:29 TU"
$ws.Range("O8").Value = "8.4 °C"

# Row 9
$ws.Range("E9").Value = "2026-02-05 18:48:26"

# Row 10
$ws.Range("E10").Value = "2026-02-05 18:48:28"
$ws.Range("O10").Value = "7.7 °C"

# Row 11
$ws.Range("E11").Value = "2026-02-05 18:48:31"
$ws.Range("J11").Value = "995.1 hPa"
$ws.Range("L11").Value = "39.2 km/h - 264º 18:11 TU"
$ws.Range("M11").Value = "4.9 °C 18:12 TU"
$ws.Range("O11").Value = "0.3 °C"

# Row 12
$ws.Range("E12").Value = "2026-02-05 18:48:34"
$ws.Range("H12").Value = "90%"
$ws.Range("L12").Value = "46.1 km/h - 285º 18:08 TU"
$ws.Range("M12").Value = "16.6 °C 18:25 TU"
$ws.Range("O12").Value = "9.7 °C"

# Row 13
$ws.Range("E13").Value = "2026-02-05 18:48:37"
$ws.Range("O13").Value = "7.7 °C"

# Row 14
$ws.Range("E14").Value = "2026-02-05 18:48:39"
$ws.Range("H14").Value = "72%"
$ws.Range("I14").Value = "6.7 mm"

# Row 15
$ws.Range("E15").Value = "2026-02-05 18:48:42"
$ws.Range("H15").Value = "84%"
$ws.Range("J15").Value = "990.7 hPa"
$ws.Range("O15").Value = "7.9 °C"

# Row 16
$ws.Range("E16").Value = "2026-02-05 18:48:44"

# Row 17
$ws.Range("E17").Value = "2026-02-05 18:48:47"

# Row 18
$ws.Range("E18").Value = "2026-02-05 18:48:50"

# Row 19
$ws.Range("E19").Value = "2026-02-05 18:48:53"

# Row 20
$ws.Range("E20").Value = "2026-02-05 18:48:56"
$ws.Range("G20").Value = "121 cm"
$ws.Range("O20").Value = "-1.4 °C"

# Row 21
$ws.Range("E21").Value = "2026-02-05 18:48:58"
$ws.Range("H21").Value = "83%"
$ws.Range("J21").Value = "990.9 hPa"
$ws.Range("O21").Value = "6.1 °C"

# Row 22
$ws.Range("E22").Value = "2026-02-05 18:49:01"
$ws.Range("H22").Value = "90%"
$ws.Range("M22").Value = "15.3 °C 18:25 TU"
$ws.Range("O22").Value = "8.3 °C"

# Row 23
$ws.Range("E23").Value = "2026-02-05 18:49:03"
$ws.Range("J23").Value = "990.1 hPa"

# Row 24
$ws.Range("E24").Value = "2026-02-05 18:49:06"
$ws.Range("J24").Value = "989.2 hPa"
$ws.Range("O24").Value = "10.3 °C"

# Row 25
$ws.Range("E25").Value = "2026-02-05 18:49:09"
$ws.Range("J25").Value = "994.2 hPa"
$ws.Range("M25").Value = "2.6 °C 18:27 TU"
$ws.Range("O25").Value = "0.5 °C"

# Row 26
$ws.Range("E26").Value = "2026-02-05 18:49:12"
$ws.Range("O26").Value = "-0.8 °C"

# Row 27
$ws.Range("E27").Value = "2026-02-05 18:49:15"
$ws.Range("J27").Value = "990.4 hPa"
$ws.Range("O27").Value = "8.6 °C"

# Row 28
$ws.Range("E28").Value = "2026-02-05 18:49:17"
$ws.Range("J28").Value = "993.3 hPa"
$ws.Range("O28").Value = "2.3 °C"

# Row 29
$ws.Range("E29").Value = "2026-02-05 18:49:20"
$ws.Range("O29").Value = "8.8 °C"

# Row 30
$ws.Range("E30").Value = "2026-02-05 18:49:23"
$ws.Range("H30").Value = "67%"
$ws.Range("K30").Value = "1.3 MJ/m2"

# Row 31
$ws.Range("E31").Value = "2026-02-05 18:49:25"
$ws.Range("I31").Value = "18.6 mm"

# Row 32
$ws.Range("E32").Value = "2026-02-05 18:49:28"
$ws.Range("H32").Value = "81%"
$ws.Range("J32").Value = "991.8 hPa"
$ws.Range("O32").Value = "11.9 °C"

# Row 33
$ws.Range("E33").Value = "2026-02-05 18:49:31"
$ws.Range("H33").Value = "88%"
$ws.Range("O33").Value = "8.9 °C"

# Row 34
$ws.Range("E34").Value = "2026-02-05 18:49:33"
$ws.Range("H34").Value = "97%"
$ws.Range("O34").Value = "3.6 °C"

# Row 35
$ws.Range("E35").Value = "2026-02-05 18:49:35"
$ws.Range("I35").Value = "4.6 mm"

# Row 36
$ws.Range("E36").Value = "2026-02-05 18:49:38"
$ws.Range("O36").Value = "10.1 °C"
